$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 13
$ws.Range("D3").Value = "Dinafex 60mg Tablet"

$ws.Range("C4").Value = 14
$ws.Range("D4").Value = "Dinafex 120mg Tablet"

$ws.Range("D5").Value = "Dinafex 180mg Tablet"

$ws.Range("C7").Value = 19

$ws.Range("C8").Value = 22
$ws.Range("D8").Value = "Etorix 120mg Tablet"
$ws.Range("E8").Value = "20's"

$ws.Range("D9").Value = "Etorix 90mg Tablet"
$ws.Range("E9").Value = "30's"

$ws.Range("C10").Value = 27

$ws.Range("C11").Value = 28
$ws.Range("D11").Value = "Flucloxin 500mg Capsule - 36's"
$ws.Range("E11").Value = "36 's"

$ws.Range("C12").Value = 29
$ws.Range("D12").Value = "Flucloxin 500mg Capsule"
$ws.Range("E12").Value = "30 's"

$ws.Range("D15").Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Range("E15").Value = "4's"

$ws.Range("C16").Value = 40
$ws.Range("D16").Value = "Ketonic 10mg Tablet"
$ws.Range("E16").Value = "20's"

$ws.Range("D17").Value = "Kynol D 25mg Tablet"
$ws.Range("E17").Value = "60 's"

$ws.Range("D18").Value = "Kynol TR 200mg Capsule"
$ws.Range("E18").Value = "30 's"

$ws.Range("E20").Value = "30 's"

$ws.Range("C21").Value = 63

$ws.Range("C22").Value = 67

$ws.Range("C23").Value = 73

$ws.Range("C25").Value = 102

$ws.Range("C26").Value = 103

$ws.Range("C27").Value = 104
$ws.Range("D27").Value = "Zithrox 30ml Dry Suspension"
$ws.Range("E27").Value = "30ml"

$ws.Range("D28").Value = "Zithrox 250mg Tablet - 6's"
$ws.Range("E28").Value = "6's"
